$wb = $excel.ActiveWorkbook

# "Latest Handoff Datetime" column (D) for the row whose "Latest Handoff File" (C)
# is the 49c10215-....xlf file has a new handoff recorded, updating the timestamp.

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-01-19 05:25:27"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-01-19 05:25:37"
